$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$row = 10

$ws.Cells.Item($row, 1).Value = "2025-08-28T06:34:14.306526"
$ws.Cells.Item($row, 2).Value = 8
$ws.Cells.Item($row, 3).Value = "全案件リスト"
$ws.Cells.Item($row, 4).Value = 62.5
$ws.Cells.Item($row, 5).Value = 3
$ws.Cells.Item($row, 6).Value = 2
$ws.Cells.Item($row, 7).Value = 8
